$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19, shifting existing rows 19-26 down to 20-27.
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the new price record (same Mercado/Region/
# Categoria/Origen/Clasificacion as the surrounding rows).
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44806
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 100112013
$ws.Range("G19").Value = "Alcachofa"
$ws.Range("H19").Value = "Madrigal"
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 80
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 13000
$ws.Range("N19").Value = "$/caja 50 unidades"
$ws.Range("O19").Value = "Provincia de Limarí"
$ws.Range("P19").Value = 260
$ws.Range("Q19").Value = 50
$ws.Range("R19").Value = "Hortaliza"
